$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# Row 2 - anchor formulas (not part of a shared group)
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("H2").Formula = "=SUM(G2:G11)"
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# Row 3 - standalone formula
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"

# Rows 4-15 - shared formula group
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# Update the selection/view to match the edited workbook
$ws.Range("B1").Select() | Out-Null
$ws.Range("J2:K2").Select() | Out-Null
